$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3 into the new row 4 (preserves cell types/values, incl. the
# empty "inlineStr" placeholder cells in P/R/T/U) before touching row 3.
$ws.Range("A3:U3").Copy($ws.Range("A4:U4"))

# Row 4 reflects the next scrape: new timestamp and a different market status.
$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("D4").Value = "未开盘"

# Row 3 loses its trailing empty placeholder cells (P3/R3/T3/U3).
$ws.Range("P3").ClearContents()
$ws.Range("R3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("U3").ClearContents()
